# Default timezone fix in tests.
#
# 1. "Game" sheet: the start-time cell (B3) used to hold a literal date
#    serial number (2020-07-01). It becomes a plain text value carrying an
#    explicit UTC offset instead, so downstream readers stop guessing the
#    local timezone.
# 2. "Level 1" sheet: materialize the blank separator row (row 5) between
#    the level header rows and the code rows, mirroring "Level 2" which
#    already has it.
# 3. Refocus the workbook back on the "Game" sheet/cell, which is what was
#    selected when the fixture was last saved.

$wb = $excel.ActiveWorkbook

$wsGame = $wb.Worksheets.Item("Game")
$wsLevel1 = $wb.Worksheets.Item("Level 1")

# 1. Replace the numeric date serial with an explicit textual timestamp.
$wsGame.Range("B3").Value = "2020-07-01 00:00:00 UTC+3"

# 2. Make row 5 on "Level 1" an explicit (blank) row, like "Level 2" has.
$wsLevel1.Range("A5:B5").NumberFormat = "General"

# 3. Restore the active sheet/selection state.
$wsGame.Activate() | Out-Null
$wsGame.Range("B4").Select() | Out-Null
